# Fill in the third reviewer's (Tiago Ventura) scores for each of the
# four songs (columns D, G, J, M) for all 20 evaluation rows (3-22).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$D = @(1,3,4,2,1,2,2,5,2,4,2,2,3,2,4,2,2,4,2,3)
$G = @(3,3,5,4,2,5,5,2,1,4,2,2,3,2,4,5,1,4,2,4)
$J = @(4,2,2,1,3,3,3,1,2,2,4,3,1,1,4,3,5,5,3,1)
$M = @(5,5,5,5,5,2,5,4,5,1,5,4,3,5,1,2,5,5,5,4)

for ($i = 0; $i -lt $D.Length; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 4).Value = $D[$i]
    $ws.Cells.Item($row, 7).Value = $G[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
    $ws.Cells.Item($row, 13).Value = $M[$i]
}

# Reflect the author's final cell selection from the edit session.
$null = $ws.Range("L27").Select()
